# Applies numeric updates to the Leve profit-tracking sheets (H:N columns)
# as refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ALC!74 (Leve Item ID 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3124.6667
$ws.Range("I74").Value = 3108.182
$ws.Range("J74").Value = 3170
$ws.Range("K74").Value = 3108.182
$ws.Range("L74").Value = 3170
$ws.Range("M74").Value = -2172.182
$ws.Range("N74").Value = -5042

# ALC!77 (Leve Item ID 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3124.6667
$ws.Range("I77").Value = 3108.182
$ws.Range("J77").Value = 3170
$ws.Range("K77").Value = 15540.91
$ws.Range("L77").Value = 15850
$ws.Range("M77").Value = -10860.91
$ws.Range("N77").Value = -25210

# ALC!113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1804.7778
$ws.Range("I113").Value = 1825.4445
$ws.Range("J113").Value = 1784.1111
$ws.Range("K113").Value = 1825.4445
$ws.Range("L113").Value = 1784.1111
$ws.Range("M113").Value = 1428.5555
$ws.Range("N113").Value = -8292.1111

# ALC!135 (Leve Item ID 44047)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1325.72
$ws.Range("I135").Value = 533.58826
$ws.Range("K135").Value = 4802.29434
$ws.Range("M135").Value = -2267.29434

# ARM!110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 37078580
$ws.Range("J110").Value = 893.6667
$ws.Range("L110").Value = 893.6667
$ws.Range("N110").Value = -4983.6667

# ARM!132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4501.696
$ws.Range("I132").Value = 5180.4375
$ws.Range("J132").Value = 2950.2856
$ws.Range("K132").Value = 15541.3125
$ws.Range("L132").Value = 8850.856800000001
$ws.Range("M132").Value = -13011.3125
$ws.Range("N132").Value = -13910.8568

# BSM!20 (Leve Item ID 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 85659.086
$ws.Range("I20").Value = 127675
$ws.Range("J20").Value = 1627.25
$ws.Range("K20").Value = 127675
$ws.Range("L20").Value = 1627.25
$ws.Range("M20").Value = -127428
$ws.Range("N20").Value = -2121.25

# BSM!94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 531.56525
$ws.Range("I94").Value = 382.70587
$ws.Range("J94").Value = 953.3333
$ws.Range("K94").Value = 382.70587
$ws.Range("L94").Value = 953.3333
$ws.Range("M94").Value = 68.29413
$ws.Range("N94").Value = -1855.3333

# CRP!22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 256.81818
$ws.Range("I22").Value = 182.5
$ws.Range("K22").Value = 182.5
$ws.Range("M22").Value = 167.5

# CRP!86 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3920.5
$ws.Range("I86").Value = 3575
$ws.Range("J86").Value = 4150.8335
$ws.Range("K86").Value = 3575
$ws.Range("L86").Value = 4150.8335
$ws.Range("M86").Value = -2452
$ws.Range("N86").Value = -6396.8335

# CRP!89 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3920.5
$ws.Range("I89").Value = 3575
$ws.Range("J89").Value = 4150.8335
$ws.Range("K89").Value = 17875
$ws.Range("L89").Value = 20754.1675
$ws.Range("M89").Value = -12259
$ws.Range("N89").Value = -31986.1675

# CRP!107 (Leve Item ID 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 539
$ws.Range("I107").Value = 548.65
$ws.Range("J107").Value = 522.9167
$ws.Range("K107").Value = 548.65
$ws.Range("L107").Value = 522.9167
$ws.Range("M107").Value = 1371.35
$ws.Range("N107").Value = -4362.9167

# CRP!118 (Leve Item ID 26112)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 49863.332
$ws.Range("J118").Value = 49863.332
$ws.Range("L118").Value = 49863.332
$ws.Range("N118").Value = -53177.332

# CRP!132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4183.8184
$ws.Range("I132").Value = 3950.1428
$ws.Range("J132").Value = 4592.75
$ws.Range("K132").Value = 11850.4284
$ws.Range("L132").Value = 13778.25
$ws.Range("M132").Value = -9320.428400000001
$ws.Range("N132").Value = -18838.25

# CUL!12 (Leve Item ID 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.105263
$ws.Range("I12").Value = 2.6
$ws.Range("K12").Value = 7.800000000000001
$ws.Range("M12").Value = 165.2

# CUL!34 (Leve Item ID 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1161.875
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 1299.2858
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 3897.8574
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -4065.8574

# CUL!98 (Leve Item ID 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 125355.375
$ws.Range("J98").Value = 200232.6
$ws.Range("L98").Value = 600697.8
$ws.Range("N98").Value = -603693.8

# CUL!131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 791.47473
$ws.Range("J131").Value = 821.65955
$ws.Range("L131").Value = 2464.97865
$ws.Range("N131").Value = -12544.97865

# GSM!70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50655.227
$ws.Range("I70").Value = 88320.164
$ws.Range("J70").Value = 5457.3
$ws.Range("K70").Value = 88320.164
$ws.Range("L70").Value = 5457.3
$ws.Range("M70").Value = -88050.164
$ws.Range("N70").Value = -5997.3

# GSM!73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 50655.227
$ws.Range("I73").Value = 88320.164
$ws.Range("J73").Value = 5457.3
$ws.Range("K73").Value = 88320.164
$ws.Range("L73").Value = 5457.3
$ws.Range("M73").Value = -87384.164
$ws.Range("N73").Value = -7329.3

# GSM!102 (Leve Item ID 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3310.3333
$ws.Range("I102").Value = 3010
$ws.Range("J102").Value = 3419.5454
$ws.Range("K102").Value = 3010
$ws.Range("L102").Value = 3419.5454
$ws.Range("M102").Value = -1388
$ws.Range("N102").Value = -6663.5454

# GSM!104 (Leve Item ID 18666)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 40671
$ws.Range("J104").Value = 40671
$ws.Range("L104").Value = 40671
$ws.Range("N104").Value = -47659

# GSM!113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1026.6471
$ws.Range("I113").Value = 576.4
$ws.Range("J113").Value = 1669.8572
$ws.Range("K113").Value = 576.4
$ws.Range("L113").Value = 1669.8572
$ws.Range("M113").Value = 1593.6
$ws.Range("N113").Value = -6009.8572

# GSM!118 (Leve Item ID 26172)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 19225
$ws.Range("J118").Value = 19225
$ws.Range("L118").Value = 19225
$ws.Range("N118").Value = -22539

# GSM!122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 716.25
$ws.Range("I122").Value = 417.22223
$ws.Range("J122").Value = 1254.5
$ws.Range("K122").Value = 1251.66669
$ws.Range("L122").Value = 3763.5
$ws.Range("M122").Value = 1198.33331
$ws.Range("N122").Value = -8663.5

# GSM!126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2916.7144
$ws.Range("I126").Value = 4015.3333
$ws.Range("J126").Value = 2092.75
$ws.Range("K126").Value = 12045.9999
$ws.Range("L126").Value = 6278.25
$ws.Range("M126").Value = -9575.999899999999
$ws.Range("N126").Value = -11218.25

# LTW!40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 93045.91
$ws.Range("I40").Value = 334666.66
$ws.Range("J40").Value = 2438.125
$ws.Range("K40").Value = 334666.66
$ws.Range("L40").Value = 2438.125
$ws.Range("M40").Value = -334530.66
$ws.Range("N40").Value = -2710.125

# LTW!46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2024799.6
$ws.Range("I46").Value = 494
$ws.Range("J46").Value = 3374336.8
$ws.Range("K46").Value = 494
$ws.Range("L46").Value = 3374336.8
$ws.Range("M46").Value = -306
$ws.Range("N46").Value = -3374712.8

# LTW!122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2498.4546
$ws.Range("I122").Value = 2419.2632
$ws.Range("K122").Value = 7257.7896
$ws.Range("M122").Value = -4807.7896

# WVR!112 (Leve Item ID 25836)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# WVR!122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1266.44
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 1600.4546
$ws.Range("K122").Value = 3012
$ws.Range("L122").Value = 4801.3638
$ws.Range("M122").Value = -562
$ws.Range("N122").Value = -9701.363799999999

# WVR!132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1831.2363
$ws.Range("I132").Value = 1784.093
$ws.Range("K132").Value = 5352.279
$ws.Range("M132").Value = -2822.279
